# Append a new log row (row 37) to the Nalco run log, mirroring the
# formatting of the previous row (row 36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 36
$newRow = 37

# Copy formatting (style) from the last existing row onto the new row.
$ws.Range("A$lastRow`:H$lastRow").Copy()
$ws.Range("A$newRow`:H$newRow").PasteSpecial(-4122)

# Populate the new row's values.
$ws.Cells.Item($newRow, 1).Value = "2025-08-20 13:05:44 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-20 18:35:44 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""

Write-Output "Appended row $newRow to $($ws.Name)"
